$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value that LOOKS like a number (e.g. "0.93151") into a cell
# as literal TEXT (shared string), matching how the workbook already stores
# its "local"/"public lb" score columns (H and I) as text rather than numbers.
# We do this by writing a text-formula into a scratch cell, copying it, and
# pasting-special "Values only" into the destination - this keeps the
# destination cell's existing style untouched (no extra number-format style
# gets created, unlike NumberFormat="@" or a leading apostrophe).
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddr, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($rangeAddr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $excel.CutCopyMode = 0
    $scratch.Clear() | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Insert a new row after row 26, pushing rows 27-29 down to 28-30.
#    Excel (and this runtime) copies the formatting of the row above into the
#    freshly inserted row, which already gives the new row F/H/I columns the
#    same per-column styles used by row 26 (wrapText feature column, plain
#    score columns) - exactly what the target file expects.
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Insert()

# ---------------------------------------------------------------------------
# 2. Update existing row 26: new feature description + new local/public
#    scores (row grows taller to fit the longer wrapped text).
# ---------------------------------------------------------------------------
$ws.Range("F26").Value = "bag of sites; CountVectorizer (1, 2)-grams max-k; morning, day, evening, night, weekday, duration, year_month, is_monday, is_wednesday, is_sunday"
Set-TextValue "H26" "0.93151"
Set-TextValue "I26" "0.94249"
$ws.Rows.Item(26).RowHeight = 60

# ---------------------------------------------------------------------------
# 3. Populate the new row 27 with the next experiment's results.
# ---------------------------------------------------------------------------
$ws.Range("E27").Value = "sgdclassifier"
$ws.Range("F27").Value = "bag of sites; CountVectorizer (1, 2)-grams max-k; morning, day, evening,  duration, year, is_monday, is_wednesday, is_sunday"
Set-TextValue "H27" "0.92982"
Set-TextValue "I27" "0.94322"
$ws.Rows.Item(27).RowHeight = 45

# ---------------------------------------------------------------------------
# 4. Cosmetic: rename the logo/picture shape (image id/name bookkeeping was
#    touched in the original edit too).
# ---------------------------------------------------------------------------
if ($ws.Shapes.Count -ge 1) {
    $ws.Shapes.Item(1).Name = "Рисунок 1"
}

# ---------------------------------------------------------------------------
# 5. Update the sheet's active selection to match the author's final cursor
#    position after finishing the edits.
# ---------------------------------------------------------------------------
$ws.Range("G27").Select() | Out-Null
